$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
# for data rows 2-10. The underlying records were reshuffled across
# rows (weekly re-sort), while all other columns stay the same.

$rows = @{
    2  = @{ D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    3  = @{ D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 }
    4  = @{ D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 }
    5  = @{ D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 }
    6  = @{ D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    7  = @{ D = 44517; M = 400; N = 5500;  O = 6000;  P = 5750;  S = 2875 }
    8  = @{ D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 }
    9  = @{ D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 }
    10 = @{ D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("S$r").Value = $vals.S
}
